# Add a second worksheet "SecondTest", populate it, and make it the active tab.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginTest")

# Add the new worksheet directly after the existing one.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "SecondTest"

# Populate the new sheet's data.
$ws2.Range("B2").Value = "dataSheet"
$ws2.Range("C2").Value = "firstTable"

# Select C3 on the new sheet, matching the target selection.
$ws2.Range("C3").Select()

# Make the new sheet the active tab.
$ws2.Activate()

$wb.Save()
